$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Bitcoin"
$ws.Range("C2").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D2").Value = "'29.240.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.09%  "

$ws.Range("B3").Value = "Ethereum"
$ws.Range("C3").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D3").Value = "'1.890.51"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.72%  "

$ws.Range("B4").Value = "TetherUSD"
$ws.Range("C4").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'322.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.04%  "

$ws.Range("B6").Value = "USDC"
$ws.Range("C6").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.18%  "

$ws.Range("B7").Value = "XRP"
$ws.Range("C7").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D7").Value = "'0.4727"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.40%  "

$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "'0.4033"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.29%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.08015"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'0.9949"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.92%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'22.96"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.60%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.930.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.21%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'5.911"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.59%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'7.028"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'89.24"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "

$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").Value = "'0.06636"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.90%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.00001027"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.28%  "

$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'17.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "'29.253.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'5.484"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.38%  "

$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'11.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.64%  "

$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").Value = "'2.176"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.86%  "

$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "'2.124.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.27%  "

$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'154.99"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.41%  "

$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'19.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.72%  "

$ws.Range("B28").Value = "InternetComputer(DFINITY)"
$ws.Range("C28").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D28").Value = "'5.938"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.02%  "

$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "'2.078"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.62%  "

$ws.Range("B30").Value = "BitcoinCash"
$ws.Range("C30").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D30").Value = "'117.08"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").Value = "'1.021"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.09414"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "

$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").Value = "'3.533"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.34%  "

$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "'1.374"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.55%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "'5.349"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "

$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02241"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.46%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "'0.06038"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.08%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.170"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.61%  "

$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'7.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.31%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.5816"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.85%  "

$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.1826"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "'10.03"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.29%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'1.286"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.24%  "

$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.07703"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "'2.366"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.88%  "

$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'12.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.5469"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.906"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "

$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'113.23"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "'0.2965"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'43.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.82%  "

